$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "27÷5=5, 2"
$t.Cell(1,2).Range.Text = "64÷8=8, 0"
$t.Cell(1,3).Range.Text = "85÷3=28, 1"
$t.Cell(1,4).Range.Text = "87÷4=21, 3"
$t.Cell(1,5).Range.Text = "65÷3=21, 2"

$t.Cell(5,1).Range.Text = "62÷7=8, 6"
$t.Cell(5,2).Range.Text = "59÷7=8, 3"
$t.Cell(5,3).Range.Text = "17÷5=3, 2"
$t.Cell(5,4).Range.Text = "79÷5=15, 4"
$t.Cell(5,5).Range.Text = "21÷4=5, 1"

$t.Cell(9,1).Range.Text = "42÷2=21, 0"
$t.Cell(9,2).Range.Text = "42÷9=4, 6"
$t.Cell(9,3).Range.Text = "44÷7=6, 2"
$t.Cell(9,4).Range.Text = "86÷5=17, 1"
$t.Cell(9,5).Range.Text = "19÷8=2, 3"

$t.Cell(13,1).Range.Text = "78÷8=9, 6"
$t.Cell(13,2).Range.Text = "36÷8=4, 4"
$t.Cell(13,3).Range.Text = "63÷3=21, 0"
$t.Cell(13,4).Range.Text = "93÷3=31, 0"
$t.Cell(13,5).Range.Text = "69÷8=8, 5"

$t.Cell(17,1).Range.Text = "44÷7=6, 2"
$t.Cell(17,2).Range.Text = "71÷7=10, 1"
$t.Cell(17,3).Range.Text = "59÷9=6, 5"
$t.Cell(17,4).Range.Text = "37÷5=7, 2"
$t.Cell(17,5).Range.Text = "97÷8=12, 1"

